# Update Bosnia Herzegovina Premier Liga odds data (swap duplicated/mis-ordered
# match rows so each row's betting data matches its correct id/teams).
# Commit: Atualizacao de bases das ligas, do dia: 06-04-2024 as 01:36
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)


# Rows 49 and 50 had their betting-data columns (B:AC) swapped;
# rank column A is left untouched.
$ws.Range("B49").Value = 6865311
$ws.Range("F49").Value = 'Sloga'
$ws.Range("G49").Value = 'GOSK Gabela'
$ws.Range("H49").Value = 3
$ws.Range("I49").Value = 2
$ws.Range("J49").Value = 'H'
$ws.Range("K49").Value = 1.833
$ws.Range("L49").Value = 3.6
$ws.Range("M49").Value = 3.4
$ws.Range("N49").Value = 1.909
$ws.Range("O49").Value = 3.4
$ws.Range("P49").Value = 3.3
$ws.Range("Q49").Value = -0.5
$ws.Range("R49").Value = 1.925
$ws.Range("S49").Value = 1.875
$ws.Range("T49").Value = 2.25
$ws.Range("U49").Value = 1.825
$ws.Range("V49").Value = 1.975
$ws.Range("W49").Value = 0.909
$ws.Range("X49").Value = -1
$ws.Range("Y49").Value = -1
$ws.Range("Z49").Value = 0.925
$ws.Range("AA49").Value = -1
$ws.Range("AB49").Value = 0.825
$ws.Range("AC49").Value = -1

$ws.Range("B50").Value = 6865310
$ws.Range("F50").Value = 'NK Igman Konjic'
$ws.Range("G50").Value = 'Zrinjski Mostar'
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 2
$ws.Range("J50").Value = 'A'
$ws.Range("K50").Value = 3.4
$ws.Range("L50").Value = 3.6
$ws.Range("M50").Value = 1.833
$ws.Range("N50").Value = 4.75
$ws.Range("O50").Value = 4.75
$ws.Range("P50").Value = 1.45
$ws.Range("Q50").Value = 1.25
$ws.Range("R50").Value = 1.775
$ws.Range("S50").Value = 2.025
$ws.Range("T50").Value = 2.75
$ws.Range("U50").Value = 1.85
$ws.Range("V50").Value = 1.95
$ws.Range("W50").Value = -1
$ws.Range("X50").Value = -1
$ws.Range("Y50").Value = 0.45
$ws.Range("Z50").Value = -1
$ws.Range("AA50").Value = 1.025
$ws.Range("AB50").Value = -1
$ws.Range("AC50").Value = 0.95

# Rows 76 and 77 had their betting-data columns (B:AC) swapped;
# rank column A is left untouched.
$ws.Range("B76").Value = 6865328
$ws.Range("F76").Value = 'Siroki Brijeg'
$ws.Range("G76").Value = 'NK Posusje'
$ws.Range("H76").Value = 1
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 'D'
$ws.Range("K76").Value = 2
$ws.Range("L76").Value = 3
$ws.Range("M76").Value = 3.5
$ws.Range("N76").Value = 2.1
$ws.Range("O76").Value = 3
$ws.Range("P76").Value = 3.3
$ws.Range("Q76").Value = -0.25
$ws.Range("R76").Value = 1.825
$ws.Range("S76").Value = 1.975
$ws.Range("T76").Value = 2
$ws.Range("U76").Value = 1.825
$ws.Range("V76").Value = 1.975
$ws.Range("W76").Value = -1
$ws.Range("X76").Value = 2
$ws.Range("Y76").Value = -1
$ws.Range("Z76").Value = -0.5
$ws.Range("AA76").Value = 0.4875
$ws.Range("AB76").Value = 0
$ws.Range("AC76").Value = 0

$ws.Range("B77").Value = 6865377
$ws.Range("F77").Value = 'Zrinjski Mostar'
$ws.Range("G77").Value = 'FK Tuzla City'
$ws.Range("H77").Value = 3
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 'H'
$ws.Range("K77").Value = 1.333
$ws.Range("L77").Value = 5
$ws.Range("M77").Value = 6
$ws.Range("N77").Value = 1.166
$ws.Range("O77").Value = 6.5
$ws.Range("P77").Value = 13
$ws.Range("Q77").Value = -2
$ws.Range("R77").Value = 1.9
$ws.Range("S77").Value = 1.9
$ws.Range("T77").Value = 3.25
$ws.Range("U77").Value = 1.95
$ws.Range("V77").Value = 1.85
$ws.Range("W77").Value = 0.1659999999999999
$ws.Range("X77").Value = -1
$ws.Range("Y77").Value = -1
$ws.Range("Z77").Value = 0
$ws.Range("AA77").Value = 0
$ws.Range("AB77").Value = 0.95
$ws.Range("AC77").Value = -1

# Rows 87 and 88 had their betting-data columns (B:AC) swapped;
# rank column A is left untouched.
$ws.Range("B87").Value = 7505497
$ws.Range("F87").Value = 'Zeljeznicar'
$ws.Range("G87").Value = 'NK Posusje'
$ws.Range("H87").Value = 1
$ws.Range("I87").Value = 1
$ws.Range("J87").Value = 'D'
$ws.Range("K87").Value = 1.65
$ws.Range("L87").Value = 3.4
$ws.Range("M87").Value = 4.75
$ws.Range("N87").Value = 1.8
$ws.Range("O87").Value = 3.2
$ws.Range("P87").Value = 4.2
$ws.Range("Q87").Value = -0.5
$ws.Range("R87").Value = 1.825
$ws.Range("S87").Value = 1.975
$ws.Range("T87").Value = 2
$ws.Range("U87").Value = 1.75
$ws.Range("V87").Value = 2.05
$ws.Range("W87").Value = -1
$ws.Range("X87").Value = 2.2
$ws.Range("Y87").Value = -1
$ws.Range("Z87").Value = -1
$ws.Range("AA87").Value = 0.9750000000000001
$ws.Range("AB87").Value = 0
$ws.Range("AC87").Value = 0

$ws.Range("B88").Value = 7505495
$ws.Range("F88").Value = 'Sloga'
$ws.Range("G88").Value = 'Zvijezda 09'
$ws.Range("H88").Value = 1
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 'H'
$ws.Range("K88").Value = 1.444
$ws.Range("L88").Value = 4.2
$ws.Range("M88").Value = 5.5
$ws.Range("N88").Value = 1.5
$ws.Range("O88").Value = 4.2
$ws.Range("P88").Value = 5.25
$ws.Range("Q88").Value = -1
$ws.Range("R88").Value = 1.8
$ws.Range("S88").Value = 2
$ws.Range("T88").Value = 2.75
$ws.Range("U88").Value = 1.775
$ws.Range("V88").Value = 2.025
$ws.Range("W88").Value = 0.5
$ws.Range("X88").Value = -1
$ws.Range("Y88").Value = -1
$ws.Range("Z88").Value = 0
$ws.Range("AA88").Value = 0
$ws.Range("AB88").Value = -1
$ws.Range("AC88").Value = 1.025

